## Item.xlsx test-data refresh for ProductDetailsPage tests
## -----------------------------------------------------------------------
## Row 2 of Sheet1 is the "Printed Dress" sample item used by the new
## ProductDetailsPage selenium test (selects an item, reads its properties
## from this sheet, and adds it to the cart). Bring the product name/size
## code in line with the item actually shown in the UI, and normalise the
## sheet's look (explicit Calibri/black font + per-column alignment that
## matches the rest of the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- stage the four distinct cell styles used on the sheet in a scratch
#     row, then copy/paste them onto the real cells one column at a time.
#     (keeps the workbook's style table identical no matter how many
#     columns share a look, instead of re-deriving the font/format on
#     every single cell.)
$scratchGeneral = $ws.Range("A10")   # font + default (General) number format
$scratchText    = $ws.Range("B10")   # font + text ("@") number format

$scratchGeneral.Font.Name = "Calibri"
$scratchGeneral.Font.Color = 0

$scratchText.Font.Name = "Calibri"
$scratchText.Font.Color = 0
$scratchText.NumberFormat = "@"

# Name / size / color columns (A, D, E) -> font only on row 1
$scratchGeneral.Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("E1").PasteSpecial(-4122)

# Quantity column (B) -> font + right aligned
$scratchGeneral.HorizontalAlignment = -4152   # xlRight
$scratchGeneral.Copy()
$ws.Range("B1:B2").PasteSpecial(-4122)

# Code column (C) -> font + text format + left aligned
$scratchText.HorizontalAlignment = -4131      # xlLeft
$scratchText.Copy()
$ws.Range("C1:C2").PasteSpecial(-4122)

# Row 2 name / size / color (A2, D2, E2) -> font + left aligned
$scratchGeneral.HorizontalAlignment = -4131   # xlLeft
$scratchGeneral.Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("E2").PasteSpecial(-4122)

$ws.Rows.Item(10).Delete() | Out-Null

# --- row heights (explicit 15pt rows) ------------------------------------
$ws.Rows.Item(1).RowHeight = 15
$ws.Rows.Item(2).RowHeight = 15

# --- data: row 2 now describes the "Printed Dress" item ------------------
$ws.Range("A2").Value = "Printed Dress"
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = "S"

# --- misc sheet/page touch-ups --------------------------------------------
$ws.PageSetup.PaperSize = 1

$ws.Range("E2").Select() | Out-Null
